$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 3).Value = 9.018204081895052
$ws.Cells.Item(2, 4).Value = 7.607254736013832
$ws.Cells.Item(2, 5).Value = 13.20875538047568
$ws.Cells.Item(2, 6).Value = 41.98635104916717
$ws.Cells.Item(2, 7).Value = 3.685679427709956
$ws.Cells.Item(2, 9).Value = 32.25314471626252
$ws.Cells.Item(2, 10).Value = 10.47440925805957
$ws.Cells.Item(2, 12).Value = 10.38439625577701
$ws.Cells.Item(2, 15).Value = 32.46513917364776
$ws.Cells.Item(3, 3).Value = 9.026284143666999
$ws.Cells.Item(3, 4).Value = 7.619166327290333
$ws.Cells.Item(3, 5).Value = 13.22722596969181
$ws.Cells.Item(3, 6).Value = 41.76131505605763
$ws.Cells.Item(3, 7).Value = 3.689264724329759
$ws.Cells.Item(3, 9).Value = 32.124723341574
$ws.Cells.Item(3, 10).Value = 10.49420576174005
$ws.Cells.Item(3, 12).Value = 10.40432801276019
$ws.Cells.Item(3, 15).Value = 32.31302885946928
$ws.Cells.Item(4, 3).Value = 9.032606829518802
$ws.Cells.Item(4, 4).Value = 7.62717542983407
$ws.Cells.Item(4, 5).Value = 13.24046196189951
$ws.Cells.Item(4, 6).Value = 41.63325448983664
$ws.Cells.Item(4, 7).Value = 3.691583311878482
$ws.Cells.Item(4, 9).Value = 32.05343635976055
$ws.Cells.Item(4, 10).Value = 10.5076110907301
$ws.Cells.Item(4, 12).Value = 10.41772010662014
$ws.Cells.Item(4, 15).Value = 32.22733609144289
$ws.Cells.Item(5, 3).Value = 9.035525839060053
$ws.Cells.Item(5, 4).Value = 7.630614320095173
$ws.Cells.Item(5, 5).Value = 13.24633228767227
$ws.Cells.Item(5, 6).Value = 41.58365123480674
$ws.Cells.Item(5, 7).Value = 3.692557734000466
$ws.Cells.Item(5, 9).Value = 32.02630787008275
$ws.Cells.Item(5, 10).Value = 10.51338847736704
$ws.Cells.Item(5, 12).Value = 10.4234679692319
$ws.Cells.Item(5, 15).Value = 32.19437749993989
$ws.Cells.Item(6, 3).Value = 9.036031224016575
$ws.Cells.Item(6, 4).Value = 7.631195929597551
$ws.Cells.Item(6, 5).Value = 13.24733583348857
$ws.Cells.Item(6, 6).Value = 41.57557170401282
$ws.Cells.Item(6, 7).Value = 3.692721325708813
$ws.Cells.Item(6, 9).Value = 32.02191977967218
$ws.Cells.Item(6, 10).Value = 10.51436681616353
$ws.Cells.Item(6, 12).Value = 10.42443995059596
$ws.Cells.Item(6, 15).Value = 32.18902391359732
$ws.Cells.Item(7, 3).Value = 9.032644809586625
$ws.Cells.Item(7, 4).Value = 7.627221098511516
$ws.Cells.Item(7, 5).Value = 13.24053920161131
$ws.Cells.Item(7, 6).Value = 41.63257501756446
$ws.Cells.Item(7, 7).Value = 3.691596333369996
$ws.Cells.Item(7, 9).Value = 32.05306269155825
$ws.Cells.Item(7, 10).Value = 10.50768773245974
$ws.Cells.Item(7, 12).Value = 10.41779644777097
$ws.Cells.Item(7, 15).Value = 32.22688362618249
$ws.Cells.Item(8, 3).Value = 9.020707590515196
$ws.Cells.Item(8, 4).Value = 7.611217696648441
$ws.Cells.Item(8, 5).Value = 13.21473077934292
$ws.Cells.Item(8, 6).Value = 41.90668478266934
$ws.Cells.Item(8, 7).Value = 3.686891377877417
$ws.Cells.Item(8, 9).Value = 32.2073074732432
$ws.Cells.Item(8, 10).Value = 10.48097575740668
$ws.Cells.Item(8, 12).Value = 10.39102942067715
$ws.Cells.Item(8, 15).Value = 32.4111083723074
$ws.Cells.Item(9, 3).Value = 9.008095190907191
$ws.Cells.Item(9, 4).Value = 7.585340521830895
$ws.Cells.Item(9, 5).Value = 13.17915328946052
$ws.Cells.Item(9, 6).Value = 42.52252545928054
$ws.Cells.Item(9, 7).Value = 3.678590008459284
$ws.Cells.Item(9, 9).Value = 32.56883389625045
$ws.Cells.Item(9, 10).Value = 10.43850311359026
$ws.Cells.Item(9, 12).Value = 10.34768134080424
$ws.Cells.Item(9, 15).Value = 32.83228033051583
$ws.Cells.Item(10, 3).Value = 9.005399333771292
$ws.Cells.Item(10, 4).Value = 7.569669068581351
$ws.Cells.Item(10, 5).Value = 13.16217310825694
$ws.Cells.Item(10, 6).Value = 43.01989616275348
$ws.Cells.Item(10, 7).Value = 3.673048073627446
$ws.Cells.Item(10, 9).Value = 32.86899788218268
$ws.Cells.Item(10, 10).Value = 10.4133257084208
$ws.Cells.Item(10, 12).Value = 10.32138811905211
$ws.Cells.Item(10, 15).Value = 33.1764319951374
$ws.Cells.Item(11, 3).Value = 9.005595722851933
$ws.Cells.Item(11, 4).Value = 7.563261710337445
$ws.Cells.Item(11, 5).Value = 13.15643464630134
$ws.Cells.Item(11, 6).Value = 43.25524313501472
$ws.Cells.Item(11, 7).Value = 3.670646410026933
$ws.Cells.Item(11, 9).Value = 33.01269722577012
$ws.Cells.Item(11, 10).Value = 10.40317753376348
$ws.Cells.Item(11, 12).Value = 10.31062878041568
$ws.Cells.Item(11, 15).Value = 33.34009884886686
$ws.Cells.Item(12, 3).Value = 9.005874185416371
$ws.Cells.Item(12, 4).Value = 7.560938908359017
$ws.Cells.Item(12, 5).Value = 13.15454686261362
$ws.Cells.Item(12, 6).Value = 43.34560911248923
$ws.Cells.Item(12, 7).Value = 3.669754018713045
$ws.Cells.Item(12, 9).Value = 33.06810788966703
$ws.Cells.Item(12, 10).Value = 10.39952209487557
$ws.Cells.Item(12, 12).Value = 10.30672695497597
$ws.Cells.Item(12, 15).Value = 33.40305771350333
$ws.Cells.Item(13, 3).Value = 9.005805145158503
$ws.Cells.Item(13, 4).Value = 7.561434564993533
$ws.Cells.Item(13, 5).Value = 13.15494075009884
$ws.Cells.Item(13, 6).Value = 43.32609281546005
$ws.Cells.Item(13, 7).Value = 3.669945453721447
$ws.Cells.Item(13, 9).Value = 33.05613052324693
$ws.Cells.Item(13, 10).Value = 10.40030102565641
$ws.Cells.Item(13, 12).Value = 10.30755961536855
$ws.Cells.Item(13, 15).Value = 33.38945539862701
$ws.Cells.Item(14, 3).Value = 9.005614544916913
$ws.Cells.Item(14, 4).Value = 7.56306853864935
$ws.Cells.Item(14, 5).Value = 13.15627362253667
$ws.Cells.Item(14, 6).Value = 43.26265300682729
$ws.Cells.Item(14, 7).Value = 3.670572650996963
$ws.Cells.Item(14, 9).Value = 33.01723615636221
$ws.Cells.Item(14, 10).Value = 10.40287304306724
$ws.Cells.Item(14, 12).Value = 10.31030431932251
$ws.Cells.Item(14, 15).Value = 33.3452590760742
$ws.Cells.Item(15, 3).Value = 9.005524359730131
$ws.Cells.Item(15, 4).Value = 7.564082869342636
$ws.Cells.Item(15, 5).Value = 13.1571271821444
$ws.Cells.Item(15, 6).Value = 43.22395457002578
$ws.Cells.Item(15, 7).Value = 3.670959047002839
$ws.Cells.Item(15, 9).Value = 32.99354075721445
$ws.Cells.Item(15, 10).Value = 10.40447288324392
$ws.Cells.Item(15, 12).Value = 10.31200798617244
$ws.Cells.Item(15, 15).Value = 33.31831416865787
$ws.Cells.Item(16, 3).Value = 9.005415088211141
$ws.Cells.Item(16, 4).Value = 7.570102305024406
$ws.Cells.Item(16, 5).Value = 13.16258806578165
$ws.Cells.Item(16, 6).Value = 43.00469336455859
$ws.Cells.Item(16, 7).Value = 3.673207421731791
$ws.Cells.Item(16, 9).Value = 32.85974788798367
$ws.Cells.Item(16, 10).Value = 10.41401515957132
$ws.Cells.Item(16, 12).Value = 10.32211542192461
$ws.Cells.Item(16, 15).Value = 33.16587562754413
$ws.Cells.Item(17, 3).Value = 9.005712128021628
$ws.Cells.Item(17, 4).Value = 7.573979693102213
$ws.Cells.Item(17, 5).Value = 13.16644657250394
$ws.Cells.Item(17, 6).Value = 42.87246869764596
$ws.Cells.Item(17, 7).Value = 3.674617233054373
$ws.Cells.Item(17, 9).Value = 32.7794789714526
$ws.Cells.Item(17, 10).Value = 10.42020316182851
$ws.Cells.Item(17, 12).Value = 10.32862356017362
$ws.Cells.Item(17, 15).Value = 33.07415235336782
$ws.Cells.Item(18, 3).Value = 9.006016928227037
$ws.Cells.Item(18, 4).Value = 7.576277805706004
$ws.Cells.Item(18, 5).Value = 13.16885283058792
$ws.Cells.Item(18, 6).Value = 42.79727609406698
$ws.Cells.Item(18, 7).Value = 3.675439362867501
$ws.Cells.Item(18, 9).Value = 32.73398585690953
$ws.Cells.Item(18, 10).Value = 10.42388519847875
$ws.Cells.Item(18, 12).Value = 10.33247998486343
$ws.Cells.Item(18, 15).Value = 33.02206739528952
$ws.Cells.Item(19, 3).Value = 9.006143148208228
$ws.Cells.Item(19, 4).Value = 7.577067583710452
$ws.Cells.Item(19, 5).Value = 13.1696996652161
$ws.Cells.Item(19, 6).Value = 42.77196664371701
$ws.Cells.Item(19, 7).Value = 3.675719656104098
$ws.Cells.Item(19, 9).Value = 32.71869967669828
$ws.Cells.Item(19, 10).Value = 10.42515298103302
$ws.Cells.Item(19, 12).Value = 10.33380514025136
$ws.Cells.Item(19, 15).Value = 33.00454892169331
$ws.Cells.Item(20, 3).Value = 9.005666647555509
$ws.Cells.Item(20, 4).Value = 7.573559908625304
$ws.Cells.Item(20, 5).Value = 13.16601648169843
$ws.Cells.Item(20, 6).Value = 42.88645572307826
$ws.Cells.Item(20, 7).Value = 3.674465993301131
$ws.Cells.Item(20, 9).Value = 32.78795406793973
$ws.Cells.Item(20, 10).Value = 10.41953172435052
$ws.Cells.Item(20, 12).Value = 10.32791905275524
$ws.Cells.Item(20, 15).Value = 33.08384720834117
$ws.Cells.Item(21, 3).Value = 9.005664994133303
$ws.Cells.Item(21, 4).Value = 7.562585793322901
$ws.Cells.Item(21, 5).Value = 13.15587438700863
$ws.Cells.Item(21, 6).Value = 43.28125352867678
$ws.Cells.Item(21, 7).Value = 3.670387965579219
$ws.Cells.Item(21, 9).Value = 33.0286336533356
$ws.Cells.Item(21, 10).Value = 10.40211249359345
$ws.Cells.Item(21, 12).Value = 10.30949345385319
$ws.Cells.Item(21, 15).Value = 33.35821429584371
$ws.Cells.Item(22, 3).Value = 9.006853241747558
$ws.Cells.Item(22, 4).Value = 7.556016883257723
$ws.Cells.Item(22, 5).Value = 13.15090837898058
$ws.Cells.Item(22, 6).Value = 43.54650422214602
$ws.Cells.Item(22, 7).Value = 3.667822168972772
$ws.Cells.Item(22, 9).Value = 33.19171344840973
$ws.Cells.Item(22, 10).Value = 10.39182053544974
$ws.Cells.Item(22, 12).Value = 10.29845660139298
$ws.Cells.Item(22, 15).Value = 33.54323145197276
$ws.Cells.Item(23, 3).Value = 9.006110413439359
$ws.Cells.Item(23, 4).Value = 7.559467714310565
$ws.Cells.Item(23, 5).Value = 13.15340684386442
$ws.Cells.Item(23, 6).Value = 43.40429471456949
$ws.Cells.Item(23, 7).Value = 3.669182518467133
$ws.Cells.Item(23, 9).Value = 33.10415728504517
$ws.Cells.Item(23, 10).Value = 10.39721365926365
$ws.Cells.Item(23, 12).Value = 10.30425528349487
$ws.Cells.Item(23, 15).Value = 33.44397646497093
$ws.Cells.Item(24, 3).Value = 9.00568679176877
$ws.Cells.Item(24, 4).Value = 7.573749478405446
$ws.Cells.Item(24, 5).Value = 13.16621034029263
$ws.Cells.Item(24, 6).Value = 42.88012961397627
$ws.Cells.Item(24, 7).Value = 3.674534332620016
$ws.Cells.Item(24, 9).Value = 32.78412043535857
$ws.Cells.Item(24, 10).Value = 10.41983489350363
$ws.Cells.Item(24, 12).Value = 10.32823720292236
$ws.Cells.Item(24, 15).Value = 33.07946214291092
$ws.Cells.Item(25, 3).Value = 9.010352202670155
$ws.Cells.Item(25, 4).Value = 7.591753207560978
$ws.Cells.Item(25, 5).Value = 13.18716891447904
$ws.Cells.Item(25, 6).Value = 42.3478185486144
$ws.Cells.Item(25, 7).Value = 3.680737427174988
$ws.Cells.Item(25, 9).Value = 32.46484992796653
$ws.Cells.Item(25, 10).Value = 10.44893360170688
$ws.Cells.Item(25, 12).Value = 10.35843133447792
$ws.Cells.Item(25, 15).Value = 32.71210494073968
